# Auto-generated Excel COM-interop edit script
# Updates Golem_Profits profitability calculations across ARM, BSM, CRP, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 4265.857
$ws.Range("I2").Value = 4265.857
$ws.Range("K2").Value = 4265.857
$ws.Range("M2").Value = -4152.857
# Row 4
$ws.Range("H4").Value = 537.5
$ws.Range("J4").Value = 537.5
$ws.Range("L4").Value = 537.5
$ws.Range("N4").Value = -769.5
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = $null
# Row 116
$ws.Range("H116").Value = 4265.857
$ws.Range("I116").Value = 4265.857
$ws.Range("K116").Value = 4265.857
$ws.Range("M116").Value = -1971.857
# Row 132
$ws.Range("H132").Value = 3856.5715
$ws.Range("I132").Value = 2999.2
$ws.Range("K132").Value = 8997.599999999999
$ws.Range("M132").Value = -6467.599999999999
# Row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = $null

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 4265.857
$ws.Range("I3").Value = 4265.857
$ws.Range("K3").Value = 4265.857
$ws.Range("M3").Value = -4151.857
# Row 86
$ws.Range("H86").Value = 3223.8572
$ws.Range("I86").Value = 3223.8572
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3223.8572
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2100.8572
$ws.Range("N86").Value = $null
# Row 89
$ws.Range("H89").Value = 3223.8572
$ws.Range("I89").Value = 3223.8572
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 16119.286
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -10503.286
$ws.Range("N89").Value = $null
# Row 105
$ws.Range("H105").Value = 1288.8
# Row 107
$ws.Range("H107").Value = 101050
$ws.Range("I107").Value = 134233.33
$ws.Range("K107").Value = 134233.33
$ws.Range("M107").Value = -132313.33
# Row 134
$ws.Range("H134").Value = 2856.875
$ws.Range("I134").Value = 2836.8572
$ws.Range("K134").Value = 8510.5716
$ws.Range("M134").Value = -5975.571599999999

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 233.66667
$ws.Range("I7").Value = 42.666668
$ws.Range("K7").Value = 42.666668
$ws.Range("M7").Value = 70.333332
# Row 39
$ws.Range("H39").Value = 3000
$ws.Range("I39").Value = 3000
$ws.Range("K39").Value = 3000
$ws.Range("M39").Value = -2609
# Row 49
$ws.Range("H49").Value = 3000
$ws.Range("I49").Value = 3000
$ws.Range("K49").Value = 3000
$ws.Range("M49").Value = -2818
# Row 86
$ws.Range("H86").Value = 166669090
$ws.Range("I86").Value = 250002750
$ws.Range("J86").Value = 1749.5
$ws.Range("K86").Value = 250002750
$ws.Range("L86").Value = 1749.5
$ws.Range("M86").Value = -250001627
$ws.Range("N86").Value = -3995.5
# Row 89
$ws.Range("H89").Value = 166669090
$ws.Range("I89").Value = 250002750
$ws.Range("J89").Value = 1749.5
$ws.Range("K89").Value = 1250013750
$ws.Range("L89").Value = 8747.5
$ws.Range("M89").Value = -1250008134
$ws.Range("N89").Value = -19979.5

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 166673330
# Row 73
$ws.Range("H73").Value = 166673330
# Row 97
$ws.Range("H97").Value = 5000
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = $null
# Row 107
$ws.Range("H107").Value = 47620564
$ws.Range("I107").Value = 200.33333
$ws.Range("K107").Value = 200.33333
$ws.Range("M107").Value = 1719.66667
# Row 113
$ws.Range("H113").Value = 2950
$ws.Range("I113").Value = 2950
$ws.Range("K113").Value = 2950
$ws.Range("M113").Value = -780
# Row 122
$ws.Range("H122").Value = 2809.3333
$ws.Range("I122").Value = 2748.5
$ws.Range("K122").Value = 8245.5
$ws.Range("M122").Value = -5795.5
# Row 132
$ws.Range("H132").Value = 4125.625
$ws.Range("I132").Value = 4125.625
$ws.Range("K132").Value = 12376.875
$ws.Range("M132").Value = -9846.875

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2232.6667
$ws.Range("I40").Value = 1350
$ws.Range("J40").Value = 3998
$ws.Range("K40").Value = 1350
$ws.Range("L40").Value = 3998
$ws.Range("M40").Value = -1214
$ws.Range("N40").Value = -4270
# Row 68
$ws.Range("H68").Value = 4411.643
$ws.Range("I68").Value = 4744.222
$ws.Range("J68").Value = 3813
$ws.Range("K68").Value = 4744.222
$ws.Range("L68").Value = 3813
$ws.Range("M68").Value = -3995.222
$ws.Range("N68").Value = -5311
# Row 71
$ws.Range("H71").Value = 4411.643
$ws.Range("I71").Value = 4744.222
$ws.Range("J71").Value = 3813
$ws.Range("K71").Value = 23721.11
$ws.Range("L71").Value = 19065
$ws.Range("M71").Value = -19977.11
$ws.Range("N71").Value = -26553
# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = $null
# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = $null
# Row 93
$ws.Range("H93").Value = 20838972
$ws.Range("I93").Value = 27782988
$ws.Range("K93").Value = 27782988
$ws.Range("M93").Value = -27781740
# Row 100
$ws.Range("H100").Value = 5000
$ws.Range("I100").Value = 4500
$ws.Range("K100").Value = 4500
$ws.Range("M100").Value = -3959

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 6748.5
$ws.Range("I62").Value = 5999.6665
$ws.Range("K62").Value = 5999.6665
$ws.Range("M62").Value = -5375.6665
# Row 65
$ws.Range("H65").Value = 6748.5
$ws.Range("I65").Value = 5999.6665
$ws.Range("K65").Value = 29998.3325
$ws.Range("M65").Value = -26878.3325
# Row 81
$ws.Range("H81").Value = 12001
$ws.Range("I81").Value = 12001
$ws.Range("K81").Value = 24002
$ws.Range("M81").Value = -22941
# Row 84
$ws.Range("H84").Value = 12001
$ws.Range("I84").Value = 12001
$ws.Range("K84").Value = 120010
$ws.Range("M84").Value = -114706
# Row 96
$ws.Range("H96").Value = 4997.5
$ws.Range("J96").Value = 4995
$ws.Range("L96").Value = 4995
$ws.Range("N96").Value = -7741
# Row 132
$ws.Range("H132").Value = 2974.6667
$ws.Range("I132").Value = 1200
$ws.Range("J132").Value = 3862
$ws.Range("K132").Value = 3600
$ws.Range("L132").Value = 11586
$ws.Range("M132").Value = -1070
$ws.Range("N132").Value = -16646
# Row 136
$ws.Range("H136").Value = 3567.3333
$ws.Range("I136").Value = 3401
$ws.Range("K136").Value = 10203
$ws.Range("M136").Value = -7653
